# Inclusao custo de transporte no orcamento
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 6) down to the
# new row (row 7) so the new cells pick up the same borders/number format.
$ws.Range("A6:B6").Copy() | Out-Null
$ws.Range("A7:B7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new "Transporte" budget line.
$ws.Range("A7").Value = "Transporte"
$ws.Range("B7").Value = "R$ 64,000"

# Match Excel's natural post-edit selection (one row below the new data).
$ws.Range("B8").Select() | Out-Null
